$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45: update I45, add J45
$ws.Range("I45").Value = 0.468739117536998
$ws.Range("J45").Value = 0.2150495036779461

# Row 46: update H46, add I46
$ws.Range("H46").Value = 0.5099036351493167
$ws.Range("I46").Value = 0.24

# Row 47: update G47, add H47
$ws.Range("G47").Value = 0.5604363747513331
$ws.Range("H47").Value = 0.3087982760018804

# Row 48: update F48, add G48
$ws.Range("F48").Value = 0.5999036351493168
$ws.Range("G48").Value = 0.32

# Row 49: update E49, add F49
$ws.Range("E49").Value = 0.6299036351493167
$ws.Range("F49").Value = 0.4476495795507702

# Row 50: update D50, add E50
$ws.Range("D50").Value = 0.3603773643037867
$ws.Range("E50").Value = 0.1088966743764388

# Row 51: update C51, add D51
$ws.Range("C51").Value = 0.4107440146302961
$ws.Range("D51").Value = 0.1461563307127136

# Row 52: update B52, add C52
$ws.Range("B52").Value = 0.25708246933236
$ws.Range("C52").Value = 0.09547648014918764

# Row 53: add B53
$ws.Range("B53").Value = 0.0959495356205764
